# B1--and-B2-PowerPoint.pptx edit
#
# The authored change updates the table style applied to the financial-
# documents comparison table on slide 5 (the table whose tblPr carried
# tableStyleId {42656FF4-3AD5-4372-A98E-313D9AD1F88A}) to the built-in
# table style {416832A8-07A1-4CBB-971F-965F324B2E90}.

$p = $ppt.ActivePresentation

$targetSlideIndex = 5
$oldStyleId = "{42656FF4-3AD5-4372-A98E-313D9AD1F88A}"
$newStyleId = "{416832A8-07A1-4CBB-971F-965F324B2E90}"

$applied = $false

# Prefer the slide that actually holds the old style id, but fall back to
# the expected slide index if the style can't be located (e.g. already
# updated).
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $tbl = $shape.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
                $applied = $true
            }
        }
    }
}

if (-not $applied) {
    $slide = $p.Slides.Item($targetSlideIndex)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
            $applied = $true
        }
    }
}

Write-Host ("Table style updated: " + $applied)
